# "Fin De semana" schedule update:
# swap the 6:35-6:55am / 6:55-7:15am activities in column B (rows 4 and 5),
# carrying their banded-row shading along with the text, and leave the
# selection resting on A26 afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fin De semana")
$ws.Activate()

# Rotate B4 -> scratch -> B5 -> B4 -> scratch -> B5 so the cell contents
# (value + number format + fill/shading) genuinely swap places instead of
# only the text moving.
$scratch = $ws.Range("Z1")
$ws.Range("B4").Cut($scratch)
$ws.Range("B5").Cut($ws.Range("B4"))
$scratch.Cut($ws.Range("B5"))
$scratch.Clear()

$ws.Range("A26").Select()
